# Applies the "Meanmodel workflow and new testcases (#6)" edit to
# "List of Tests.xlsx" (Tabelle1): refreshed test-case descriptions/content,
# tightened a handful of row heights, dropped the unused red-font style, and
# moved the saved selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---------------------------------------------------------------------------
# 1. Cell text content (columns A-E). Row/column numbers (column B) already
#    hold the correct values and are left untouched.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "NO"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "expected outcome"
$ws.Range("E1").Value = "Models"

$ws.Range("C2").Value = "Check if all necessary input files exists"
$ws.Range("D2").Value = "Breaks with logfile message like this:`n25-Jul-2017 13:36:56 Start input checks`n25-Jul-2017 13:36:56 ERROR: `"Children_OralSingle_IV_Multi_withTypo.xml`" does not exist`n25-Jul-2017 13:36:56 ERROR: `"Children_OralSingle_IV_Multi_withTypo.csv `" does not exist"
$ws.Range("E2").Value = "TestExample"

$ws.Range("C3").Value = "Find inconsistencies in output definition, unkonw output path, unknown units, or wrong unit for dimensions"
$ws.Range("D3").Value = "Breaks with logfile message`nERROR: Outputpath `"Organism|PeripheralVenousBlood|Hydroxy_Itraconazole|Plasma (Peripheral Venous Blood) with Typo`" could not be found in model`nERROR: For unit `"µmol/l`", there is no common dimension with display unit `"cm`"`nERROR: unit `"typo`" for seems to be no default OSPSuite unit`nERROR: unit `"typo2`" for seems to be no default OSPSuite unit`nERROR: For unit `"µmol/l`", there is no common dimension with display unit `"h`""
$ws.Range("E3").Value = "TestExample"

$ws.Range("C4").Value = "Population simulation and PK Parameter calculation of a population with single application"
$ws.Range("D4").Value = "Two csv files are generated: `nsimulations/SingleIvBolus-Results.csv`nsimulations/SingleIvBolus-PK-Analyses.csv`nautomatic test compares the result to PK-SIM export of this model. It works for simulations, fails for PK Parameter with extrapolation ???"

$ws.Range("C5").Value = "Population simulation and PK Parameter calculation of a multi application"
$ws.Range("D5").Value = "Two csv files are generated: `nSimulations/OralSingle_IV_Multi-PK-Analyses.csv`nSimulations/OralSingle_IV_Multi-Results.csv`nautomatic test compares the result to PK-SIM export of this model. It works for simulations, fails for PK Parameter with extrapolation ???"
$ws.Range("E5").Value = "TestExample"

$ws.Range("C6").Value = "Populations simulation with studydesign.csv which sets the dose according to the dose per bodyweigth"
$ws.Range("D6").Value = "Two csv files are generated: `nsimulations/SingleIvBolus-Results.csv`nsimulations/SingleIvBolus-PK-Analyses.csv`nAutomatic test checks if the simulation results is a factor 10 highe the the PK-Sim export"

$ws.Range("C7").Value = "Populations simulation with studydesign.csv dose per surface area"
$ws.Range("D7").Value = "Two csv files are generated: `nsimulations/SingleIvBolus-Results.csv`nsimulations/SingleIvBolus-PK-Analyses.csv`nAutomatic test checks if the simulation results is a factor 10 highe the the PK-Sim export"
$ws.Range("E7").Value = "7.2_BSA_Example"

$ws.Range("C8").Value = "Read nonmen file with individual timeprofiles simulate mean model and plot data vs prediction"
$ws.Range("D8").Value = "data is converted to a mat file: tmp/PO320mg/dataTP.mat  check for random individuals, if time and dv is correctly transferred`nfigures are created: time profile, same figures as in corresponding output, `n"

$ws.Range("C9").Value = "Population VPC with data for a single population"
$ws.Range("D9").Value = "figures are created: pyhsiological, time profile and pkParameter, same figures as in corresponding output, `n"

$ws.Range("C10").Value = "Population simulation for workflowtype parallelComparison"
$ws.Range("D10").Value = "figures are created: pyhsiological, time profile and pkParameter, same figures as in corresponding output, `n"

$ws.Range("C11").Value = "Population simulation for workflowtype pediatric"
$ws.Range("D11").Value = "figures are created: pyhsiological, time profile and pkParameter, same figures as in corresponding output, `n"

# Row 12 now only keeps its (empty) formatted cell - the old "Find
# inconsitencies in output cvs" text was removed.
$ws.Range("C12").Value = ""

# ---------------------------------------------------------------------------
# 2. Row heights - a few rows grew/shrank to fit the rewritten text.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 55.2
$ws.Rows.Item(7).RowHeight = 55.2
$ws.Rows.Item(8).RowHeight = 55.2
$ws.Rows.Item(9).RowHeight = 41.4
$ws.Rows.Item(10).RowHeight = 41.4
$ws.Rows.Item(11).RowHeight = 41.4

# ---------------------------------------------------------------------------
# 3. Drop the red warning font: cells C4/C6/C7 and C8:C12 all go back to the
#    plain/automatic text colour (the red font style is no longer used
#    anywhere in the workbook).
# ---------------------------------------------------------------------------
$ws.Range("C4").Font.ColorIndex = -4105
$ws.Range("C6").Font.ColorIndex = -4105
$ws.Range("C7").Font.ColorIndex = -4105
$ws.Range("C8:C12").Font.ColorIndex = -4105

# ---------------------------------------------------------------------------
# 4. View state - clear the frozen scroll position and move the remembered
#    selection from D11 to D10.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Top = 216
$win.Height = 9084
$ws.Range("D10").Select()
